# "AddNonPriceAgreementItem_FromGlobalItems" test data sheet
#
# 1. Update the remembered selection on each of the four existing sheets.
# 2. Append a brand-new worksheet "AddNonPriceAgr_GlobalCatalog" (becomes
#    sheet5.xml / sheetId 5 / rId5) holding the new scenario's header +
#    sample row, and leave it the active / tab-selected sheet.

$wb = $excel.ActiveWorkbook

# --- sheet1 "SmartForm": selection C11 -> C19 (loses tabSelected once a
#     later sheet is activated below) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("C19").Select()

# --- sheet2 "VerifyCatalogSearch": selection G1 -> C2 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("C2").Select()

# --- sheet3 "AddPriceAgrmnt_RecentOrder": selection I7 -> H2 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()
$ws3.Range("H2").Select()

# --- sheet4 "AddPriceAgrmnt_LocalCatalog": selection J6 -> K14 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Activate()
$ws4.Range("K14").Select()

# --- new sheet5 "AddNonPriceAgr_GlobalCatalog" ---
$ws5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws5.Name = "AddNonPriceAgr_GlobalCatalog"

# header row
$ws5.Range("A1").Value = "Role"
$ws5.Range("B1").Value = "Location"
$ws5.Range("C1").Value = "Global Item"
$ws5.Range("D1").Value = "ItemDescription"
$ws5.Range("E1").Value = "Quantity"
$ws5.Range("F1").Value = "Unit of Measure"
$ws5.Range("G1").Value = "UnitPrice"
$ws5.Range("H1").Value = "UpdatedUnitPrice"
$ws5.Range("I1").Value = "UpdatedUnitofMeasure"

# sample data row
$ws5.Range("A2").Value = "REQUESTOR"
$ws5.Range("B2").Value = "XEEVA -MJ"
$ws5.Range("C2").Value = "Iphone"
$ws5.Range("D2").Value = "iPhone_d_99_4"
$ws5.Range("E2").Value = 2
$ws5.Range("F2").Value = "CU-CUBIC"
$ws5.Range("G2").Value = 2
$ws5.Range("H2").Value = 3
$ws5.Range("I2").Value = "EA-EACH"

# approximate column widths matching the authored layout
$ws5.Columns.Item(1).ColumnWidth = 13.140625
$ws5.Columns.Item(2).ColumnWidth = 13.140625
$ws5.Columns.Item(3).ColumnWidth = 14.140625
$ws5.Columns.Item(4).ColumnWidth = 17.5703125
$ws5.Columns.Item(5).ColumnWidth = 13.5703125
$ws5.Columns.Item(6).ColumnWidth = 16.85546875
$ws5.Columns.Item(7).ColumnWidth = 11.5703125
$ws5.Columns.Item(8).ColumnWidth = 19.140625
$ws5.Columns.Item(9).ColumnWidth = 22.7109375

# leave the new sheet active/tab-selected with its remembered selection
$ws5.Range("F12").Select()
